$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert 6 blank rows before row 4 (old rows 4-8 move down to 10-14) ---
$ws.Rows("4:9").Insert()

# --- Step 2: duplicate the (untouched) title block rows 1-3 into the freed rows 7-9 ---
$ws.Range("A1:O3").Copy($ws.Range("A7:O7"))

# --- Step 3: rebuild the (now-shifted) header row (row 10, old row 4) into the new row 4 ---
# Copy-with-destination preserves styles; this fills row 4 exactly like old row 4 (incl. K4=Number of Boxes)
$ws.Range("A10:O10").Copy($ws.Range("A4:O4"))
# K4 must end up blank (value removed, style s="8" kept)
$ws.Range("K4").ClearContents()

# --- Step 4: update A3 value (STORAGE -> DEFAULT). Do this before step 5 so the new shared
# string "DEFAULT" is appended to the table first, matching the original authoring order. ---
$ws.Range("A3").Value = "DEFAULT"

# --- Step 5: rebuild the "old row 5" pattern (now at row 11) into new row 5, then edit its values ---
$ws.Range("A11:O11").Copy($ws.Range("A5:O5"))
$ws.Range("A5").ClearContents()
$ws.Range("B5").Value = "/DEFAULT/DEFAULT/DEFAULT"
$ws.Range("C5").Value = "DEFAULT"
$ws.Range("D5").Value = "DEFAULT"
$ws.Range("E5").Value = "/DEFAULT/DEFAULT"
$ws.Range("F5").Value = "/DEFAULT/DEFAULT/DEFAULT"
$ws.Range("G5").Value = $false
$ws.Range("H5").ClearContents()
$ws.Range("I5").ClearContents()
$ws.Range("J5").Value = "Default"
$ws.Range("K5").ClearContents()
$ws.Rows("5:5").RowHeight = 17

# --- Step 6: row 6 must stay completely empty (it was left blank by the Insert in step 1) ---

# --- Step 7: widen column A ---
$ws.Columns("A:A").ColumnWidth = 11

# --- Step 8: reset selection to A1 ---
$ws.Range("A1").Select()
